$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ANALYSIS_UNIT")

$rng = $ws.Range("A200")
$rng.Interior.TintAndShade = 0.79998168889431442
$rng.Interior.ThemeColor = 6
$rng.Value = "t"

$rng2 = $ws.Range("A201")
$rng2.Interior.Color = 16750960
$rng2.Value = "t2"
